$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 25 (content shifted up / removed)
$ws.Rows(25).Delete()

# Row 10
$ws.Range("B10").Value = "7455355 - Robson da Silva Rocha"
$ws.Range("C10").Value = "7455355 - Robson da Silva Rocha"

# Row 13
$ws.Rows(13).RowHeight = 60
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 14
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Systems of Water Supply; Water Treatment Technologies; Water Treatment in Complete Cycle; Disinfection, Direct Ascendant Filtration; Direct Descendant Filtration, Floto-filtration; Filtration in Multiple Levels, Treatment of the Generated Waste in the Water Treatment Stations."
$ws.Range("C14").Value = "Systems of Water Supply; Water Treatment Technologies; Water Treatment in Complete Cycle; Disinfection, Direct Ascendant Filtration; Direct Descendant Filtration, Floto-filtration; Filtration in Multiple Levels, Treatment of the Generated Waste in the Water Treatment Stations."

# Row 15
$ws.Rows(15).RowHeight = 120
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2022"
$ws.Range("C15").Value = "01/01/2022"

# Row 16
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "- Types of dimensioning of water distribution networks;- Water characteristics of interest for treatment: physical, chemical and bacteriological characteristics;- Potability Standard;- Water Treatment Technologies;- Constituent Units of a Water Supply System;- Groundwater Catchment and Surface Water Catchment- Railing, sand removal, pump room;- Reservation;- Distribution networks: types of network, parts and accessories;- Full Cycle Water Treatment System;- Coagulation-flocculation and Rapid Mixing;- Decantation: conventional and high rate decantation and sludge removal system;- Filtration mechanisms, filter materials and filter bottoms, filtration hydraulics, filtration with constant rate and declining rate,- Disinfection: main disinfectants, chlorination and chlor-ammonia, pre and post-chlorination,- Treatment of waste generated at stations and reuse of recovered water"
$ws.Range("C16").Value = "- Types of dimensioning of water distribution networks;- Water characteristics of interest for treatment: physical, chemical and bacteriological characteristics;- Potability Standard;- Water Treatment Technologies;- Constituent Units of a Water Supply System;- Groundwater Catchment and Surface Water Catchment- Railing, sand removal, pump room;- Reservation;- Distribution networks: types of network, parts and accessories;- Full Cycle Water Treatment System;- Coagulation-flocculation and Rapid Mixing;- Decantation: conventional and high rate decantation and sludge removal system;- Filtration mechanisms, filter materials and filter bottoms, filtration hydraulics, filtration with constant rate and declining rate,- Disinfection: main disinfectants, chlorination and chlor-ammonia, pre and post-chlorination,- Treatment of waste generated at stations and reuse of recovered water"

# Row 17
$ws.Rows(17).RowHeight = 15
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()

# Row 18
$ws.Rows(18).RowHeight = 60
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "7455355 - Robson da Silva Rocha"
$ws.Range("C18").Value = "7455355 - Robson da Silva Rocha"

# Row 19
$ws.Range("A19").Value = "Critério:"

# Row 20
$ws.Range("A20").Value = "Norma de recuperação:"

# Row 21
$ws.Rows(21).RowHeight = 120
$ws.Range("A21").Value = "Bibliografia:"

# Row 22
$ws.Rows(22).RowHeight = 15
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()

# Row 23
$ws.Rows(23).RowHeight = 30
$ws.Range("A23").Clear()
$ws.Range("B23").Value = "LOB1212 -  Química Analítica Ambiental II  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOB1212 -  Química Analítica Ambiental II  (Requisito fraco)`n"

# Row 24
$ws.Range("B24").Value = "LOB1258 -  Hidráulica Aplicada  (Requisito fraco)`n"
$ws.Range("C24").Value = "LOB1258 -  Hidráulica Aplicada  (Requisito fraco)`n"

